$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text formatting
# (values such as "0.630" or "250.86" must not be coerced to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "41.775.59"
$ws.Cells.Item(2, 5).Value = "  -1.08%  "
$ws.Cells.Item(3, 4).Value = "2.216.95"
$ws.Cells.Item(3, 5).Value = "  -1.22%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "250.86"
$ws.Cells.Item(5, 5).Value = "  +2.28%  "
$ws.Cells.Item(6, 4).Value = "0.630"
$ws.Cells.Item(6, 5).Value = "  -0.18%  "
$ws.Cells.Item(7, 4).Value = "70.35"
$ws.Cells.Item(7, 5).Value = "  +1.78%  "
$ws.Cells.Item(9, 4).Value = "0.604"
$ws.Cells.Item(9, 5).Value = "  +9.09%  "
$ws.Cells.Item(10, 4).Value = "40.04"
$ws.Cells.Item(10, 5).Value = "  +10.24%  "
$ws.Cells.Item(11, 4).Value = "0.0964"
$ws.Cells.Item(12, 4).Value = "58.24"
$ws.Cells.Item(12, 5).Value = "  -1.70%  "
$ws.Cells.Item(13, 4).Value = "7.23"
$ws.Cells.Item(13, 5).Value = "  +7.40%  "
$ws.Cells.Item(15, 4).Value = "2.546.99"
$ws.Cells.Item(15, 5).Value = "  -1.21%  "
$ws.Cells.Item(16, 4).Value = "14.97"
$ws.Cells.Item(16, 5).Value = "  -0.27%  "
$ws.Cells.Item(17, 5).Value = "  +2.79%  "
$ws.Cells.Item(18, 4).Value = "2.219.19"
$ws.Cells.Item(18, 5).Value = "  -1.27%  "
$ws.Cells.Item(19, 4).Value = "41.797.84"
$ws.Cells.Item(19, 5).Value = "  -0.85%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0961"
$ws.Cells.Item(20, 5).Value = "  -0.86%  "
$ws.Cells.Item(21, 4).Value = "6.24"
$ws.Cells.Item(21, 5).Value = "  +0.10%  "
$ws.Cells.Item(22, 4).Value = "72.58"
$ws.Cells.Item(22, 5).Value = "  -0.82%  "
$ws.Cells.Item(23, 4).Value = "234.75"
$ws.Cells.Item(23, 5).Value = "  -0.50%  "
$ws.Cells.Item(24, 5).Value = "  +1.52%  "
$ws.Cells.Item(25, 5).Value = "  +12.37%  "
$ws.Cells.Item(26, 4).Value = "11.75"
$ws.Cells.Item(26, 5).Value = "  +17.61%  "
$ws.Cells.Item(27, 5).Value = "  -0.01%  "
$ws.Cells.Item(28, 4).Value = "2.52"
$ws.Cells.Item(28, 5).Value = "  +1.47%  "
$ws.Cells.Item(29, 5).Value = "  -1.56%  "
$ws.Cells.Item(30, 4).Value = "169.99"
$ws.Cells.Item(30, 5).Value = "  -1.58%  "
$ws.Cells.Item(31, 4).Value = "20.79"
$ws.Cells.Item(31, 5).Value = "  +1.29%  "
$ws.Cells.Item(32, 4).Value = "0.122"
$ws.Cells.Item(32, 5).Value = "  +0.10%  "
$ws.Cells.Item(33, 5).Value = "  -1.86%  "
$ws.Cells.Item(34, 4).Value = "5.52"
$ws.Cells.Item(34, 5).Value = "  +4.02%  "
$ws.Cells.Item(35, 4).Value = "0.0737"
$ws.Cells.Item(35, 5).Value = "  +2.82%  "
$ws.Cells.Item(36, 4).Value = "4.67"
$ws.Cells.Item(36, 5).Value = "  -0.84%  "
$ws.Cells.Item(37, 4).Value = "26.04"
$ws.Cells.Item(37, 5).Value = "  +14.19%  "
$ws.Cells.Item(38, 4).Value = "4.05"
$ws.Cells.Item(38, 5).Value = "  +6.97%  "
$ws.Cells.Item(39, 4).Value = "0.0307"
$ws.Cells.Item(39, 5).Value = "  +8.40%  "
$ws.Cells.Item(40, 4).Value = "2.27"
$ws.Cells.Item(40, 5).Value = "  -1.33%  "
$ws.Cells.Item(41, 4).Value = "5.90"
$ws.Cells.Item(41, 5).Value = "  -0.08%  "
$ws.Cells.Item(42, 4).Value = "12.44"
$ws.Cells.Item(42, 5).Value = "  +23.78%  "
$ws.Cells.Item(43, 4).Value = "64.93"
$ws.Cells.Item(43, 5).Value = "  -3.03%  "
$ws.Cells.Item(44, 4).Value = "0.206"
$ws.Cells.Item(44, 5).Value = "  +8.48%  "
$ws.Cells.Item(45, 4).Value = "4.80"
$ws.Cells.Item(45, 5).Value = "  -5.43%  "
$ws.Cells.Item(46, 4).Value = "8.70"
$ws.Cells.Item(46, 5).Value = "  -6.28%  "
$ws.Cells.Item(47, 5).Value = "  -0.69%  "
$ws.Cells.Item(48, 4).Value = "4.66"
$ws.Cells.Item(48, 5).Value = "  -0.60%  "
$ws.Cells.Item(49, 5).Value = "  -0.36%  "
$ws.Cells.Item(50, 4).Value = "1.17"
$ws.Cells.Item(50, 5).Value = "  +5.55%  "
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).Value = "2.38"
$ws.Cells.Item(51, 5).Value = "  +2.30%  "
